# Weekly update: insert 4 new price rows (date 2022-08-25) at the top of the
# "Alcachofa" data block (row 185), pushing the existing rows down by 4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 185 (existing rows 185.. shift down to 189..)
$ws.Range("A185:A188").EntireRow.Insert()

$newDate = Get-Date -Year 2022 -Month 8 -Day 25 -Hour 0 -Minute 0 -Second 0

$rowsData = @(
    @{ Row = 185; H = "Argentina(o)"; I = "Primera"; J = 2400; K = 9500;  L = 10000; M = 9750; N = '$/caja 50 unidades'; P = 195; Q = 50 },
    @{ Row = 186; H = "Argentina(o)"; I = "Segunda"; J = 1800; K = 8500;  L = 9000;  M = 8750; N = '$/caja 70 unidades'; P = 125; Q = 70 },
    @{ Row = 187; H = "Española";     I = "Primera"; J = 1000; K = 10000; L = 11000; M = 10500; N = '$/caja 30 unidades'; P = 350; Q = 30 },
    @{ Row = 188; H = "Madrigal";     I = "Primera"; J = 1100; K = 8500;  L = 9000;  M = 8750; N = '$/caja 40 unidades'; P = 219; Q = 40 }
)

foreach ($r in $rowsData) {
    $row = $r.Row
    $ws.Range("A$row").Value = 2
    $ws.Range("B$row").Value = "Comercializadora del Agro de Limarí"
    $ws.Range("C$row").Value = "Coquimbo"
    $ws.Range("D$row").Value = $newDate
    $ws.Range("E$row").Value = 4
    $ws.Range("F$row").Value = 100112013
    $ws.Range("G$row").Value = "Alcachofa"
    $ws.Range("H$row").Value = $r.H
    $ws.Range("I$row").Value = $r.I
    $ws.Range("J$row").Value = $r.J
    $ws.Range("K$row").Value = $r.K
    $ws.Range("L$row").Value = $r.L
    $ws.Range("M$row").Value = $r.M
    $ws.Range("N$row").Value = $r.N
    $ws.Range("O$row").Value = "Provincia de Limarí"
    $ws.Range("P$row").Value = $r.P
    $ws.Range("Q$row").Value = $r.Q
    $ws.Range("R$row").Value = "Hortaliza"
}

Write-Host "Inserted 4 rows with new weekly data at row 185"
